# Update cryptos list worksheet values (price/volume columns) per the
# scheduled data refresh. Values that look numeric in column D are
# prefixed with a leading apostrophe so Excel stores them as literal
# text (matching the source data's formatting, e.g. "615.28" / "0.0361")
# instead of silently re-casting them to floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.334.26"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "3.762.57"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'615.28"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").Value = "'179.86"
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("D7").Value = "3.761.23"
$ws.Range("E7").Value = "  -0.66%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -1.53%  "
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("D11").Value = "'6.54"
$ws.Range("E11").Value = "  +3.23%  "
$ws.Range("E12").Value = "  -2.02%  "
$ws.Range("D13").Value = "'40.16"
$ws.Range("E13").Value = "  -2.25%  "
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("D15").Value = "4.391.97"
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("D16").Value = "3.766.25"
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("D17").Value = "69.393.91"
$ws.Range("E17").Value = "  -0.89%  "
$ws.Range("E18").Value = "  -2.52%  "
$ws.Range("D19").Value = "'7.49"
$ws.Range("E19").Value = "  -1.43%  "
$ws.Range("D20").Value = "'16.46"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("D21").Value = "'499.21"
$ws.Range("E21").Value = "  -3.23%  "
$ws.Range("D22").Value = "'9.43"
$ws.Range("E22").Value = "  -1.85%  "
$ws.Range("D23").Value = "'0.724"
$ws.Range("E23").Value = "  -0.51%  "
$ws.Range("E24").Value = "  +1.67%  "
$ws.Range("E25").Value = "  -2.60%  "
$ws.Range("D26").Value = "'12.94"
$ws.Range("E26").Value = "  -3.00%  "
$ws.Range("D27").Value = "'10.95"
$ws.Range("E27").Value = "  -1.67%  "
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  +0.79%  "
$ws.Range("D31").Value = "'8.10"
$ws.Range("E31").Value = "  +3.43%  "
$ws.Range("E32").Value = "  +3.36%  "
$ws.Range("D33").Value = "'30.64"
$ws.Range("E33").Value = "  -3.34%  "
$ws.Range("E34").Value = "  -1.48%  "
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("E37").Value = "  -1.49%  "
$ws.Range("E38").Value = "  +3.17%  "
$ws.Range("E39").Value = "  +4.05%  "
$ws.Range("D40").Value = "'464.34"
$ws.Range("E40").Value = "  +9.57%  "
$ws.Range("D41").Value = "'3.08"
$ws.Range("E41").Value = "  +11.90%  "
$ws.Range("E42").Value = "  -4.97%  "
$ws.Range("B43").Value = "Arweave"
$ws.Range("C43").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D43").Value = "'45.40"
$ws.Range("E43").Value = "  +2.24%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "'49.63"
$ws.Range("E44").Value = "  -3.24%  "
$ws.Range("D45").Value = "'8.59"
$ws.Range("E45").Value = "  -2.42%  "
$ws.Range("D46").Value = "2.954.53"
$ws.Range("E46").Value = "  -3.66%  "
$ws.Range("D47").Value = "'0.0361"
$ws.Range("E47").Value = "  -0.76%  "
$ws.Range("D48").Value = "'27.41"
$ws.Range("E48").Value = "  -1.07%  "
$ws.Range("D49").Value = "'138.64"
$ws.Range("E49").Value = "  +2.60%  "
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("E51").Value = "  -1.19%  "
